$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark (bookmarkStart/bookmarkEnd) from
#    the first paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append a new paragraph after the last paragraph ("Text for students"),
#    inheriting its paragraph formatting, and set its text to "2 ".
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd
$newPara = $d.Paragraphs.Add($endRange)
$newPara.Range.Text = "2 "
